$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Collected Amount *" column (column G) entirely; this shifts
# Status, Verified, Folio No* (H, I, J) one column to the left.
$ws.Range("G1").EntireColumn.Delete()

# Selection as left by the editor after the column delete
$ws.Range("G1:G1048576").Select()
